$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - new path rows (A2 already has "about-cancer/coping/feelings")
$ws.Range("A3").Value = "/grants-training/apply-grant/development"
$ws.Range("A4").Value = "/news-events/press-releases/2018/cll-ibrutinib-trial"
$ws.Range("A5").Value = "/news-events/cancer-currents-blog"
$ws.Range("A6").Value = "news-events/cancer-currents-blog/2018/cancer-moonshot-planning-to-research"
$ws.Range("A7").Value = "news-events/press-releases/2015/ilc-2015"
$ws.Range("A8").Value = "/about-nci/budget/congressional-justification"
$ws.Range("A9").Value = "/about-nci/budget/about-annual-plan"
$ws.Range("A10").Value = "about-nci/budget/fact-book/historical-trends"

# Column B - new "type" column
$ws.Range("B1").Value = "type"
$ws.Range("B2").Value = "Article"
$ws.Range("B3").Value = "Article"
$ws.Range("B4").Value = "Press Release"
$ws.Range("B5").Value = "Blog Series"
$ws.Range("B6").Value = "Blog Post"
$ws.Range("B7").Value = "Press Release"
$ws.Range("B8").Value = "General"
$ws.Range("B9").Value = "General"
$ws.Range("B10").Value = "General"

# Widen column A to fit the longer path strings
# (36.75 is the closest COM ColumnWidth that round-trips to the saved
# OOXML width of ~37.6640625 - the API quantizes width to 1/6-character steps)
$ws.Columns.Item(1).ColumnWidth = 36.75

# Give the last three path rows a distinct cell style (applyFont)
$ws.Range("A8:A10").Font.ThemeColor = 1

# Final selection, matching the saved workbook state
$ws.Range("B11").Select() | Out-Null
